# The page footer ("Ver no Jupiter Salvar em pdf Salvar em docx" and the
# "© 2020 ... Creative Commons Attribution" copyright line), plus the blank
# separator paragraph that followed them, were dropped from the end of the
# document (a site-footer rebuild). Locate that span by its text and delete
# it as a single range so the surrounding paragraphs collapse together,
# leaving just the one blank paragraph before the trailing page-break
# paragraph.

$d = $word.ActiveDocument

$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text

    if ($startPara -eq $null -and $t -like "Ver no Jupiter*") {
        $startPara = $i
    }
    if ($t -like "*Creative Commons Attribution*") {
        # Also swallow the blank paragraph right after the copyright line.
        $endPara = $i + 1
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $start = $d.Paragraphs.Item($startPara).Range.Start
    $end = $d.Paragraphs.Item($endPara).Range.End
    $r = $d.Range($start, $end)
    $r.Delete()
}
